$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Instructions Grading" and "Write Descriptive Tests" sections as Completed
$ws.Range("F3:F7").Value = "Completed"
$ws.Range("F11:F15").Value = "Completed"

# Mark the "Assessment" section as Not Completed
$ws.Range("F18:F20").Value = "Not Completed"
